$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need the Text format
# applied first, otherwise Excel would auto-convert the typed text into a
# numeric value (e.g. "1.00" -> 1), which does not match the source data
# (these columns store prices as literal text, e.g. "0.0000337", "41.50").
$textCells = @("D4", "D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D19", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D45", "D47", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin price / volume data (and the two rank swaps) per the latest run
$ws.Range("D2").Value = "66.056.75"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "3.823.45"
$ws.Range("E3").Value = "  +8.38%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "428.16"
$ws.Range("E5").Value = "  +8.24%  "
$ws.Range("D6").Value = "131.39"
$ws.Range("E6").Value = "  +5.24%  "
$ws.Range("D7").Value = "3.817.77"
$ws.Range("E7").Value = "  +8.45%  "
$ws.Range("E8").Value = "  +3.61%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("D10").Value = "0.733"
$ws.Range("E10").Value = "  +7.05%  "
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  +3.77%  "
$ws.Range("D12").Value = "0.0000337"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "41.50"
$ws.Range("E13").Value = "  +5.67%  "
$ws.Range("D14").Value = "10.43"
$ws.Range("E14").Value = "  +12.36%  "
$ws.Range("D15").Value = "4.442.28"
$ws.Range("E15").Value = "  +9.04%  "
$ws.Range("D16").Value = "15.56"
$ws.Range("E16").Value = "  +22.54%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.138"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.839.41"
$ws.Range("E18").Value = "  +8.71%  "
$ws.Range("D19").Value = "20.10"
$ws.Range("E19").Value = "  +6.55%  "
$ws.Range("E20").Value = "  +8.11%  "
$ws.Range("D21").Value = "66.335.99"
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").Value = "416.62"
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("D23").Value = "15.12"
$ws.Range("E23").Value = "  +8.18%  "
$ws.Range("D24").Value = "85.43"
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("E25").Value = "  +8.26%  "
$ws.Range("D26").Value = "37.28"
$ws.Range("E26").Value = "  +9.60%  "
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  +14.26%  "
$ws.Range("D28").Value = "3.30"
$ws.Range("E28").Value = "  +9.78%  "
$ws.Range("D29").Value = "9.45"
$ws.Range("E29").Value = "  +37.40%  "
$ws.Range("D30").Value = "5.39"
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("D31").Value = "14.11"
$ws.Range("E31").Value = "  +18.13%  "
$ws.Range("D32").Value = "711.93"
$ws.Range("E32").Value = "  +5.11%  "
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  +12.81%  "
$ws.Range("D34").Value = "2.77"
$ws.Range("E34").Value = "  +7.66%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "5.82"
$ws.Range("E35").Value = "  +42.25%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "38.95"
$ws.Range("E37").Value = "  +5.47%  "
$ws.Range("D38").Value = "0.150"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "55.68"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").Value = "0.0471"
$ws.Range("E40").Value = "  +6.87%  "
$ws.Range("D41").Value = "0.0₃0725"
$ws.Range("E41").Value = "  +15.37%  "
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("D45").Value = "3.41"
$ws.Range("E45").Value = "  +10.00%  "
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "0.322"
$ws.Range("E47").Value = "  +16.09%  "
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  +39.97%  "
$ws.Range("E49").Value = "  +6.28%  "
$ws.Range("E50").Value = "  +5.47%  "
$ws.Range("E51").Value = "  +4.34%  "
